$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Windows.Item(1).Width = 18240
$wb.Windows.Item(1).Height = 8000

$ws.Range("C3").Value = 0.1012
$ws.Range("C4").Value = 0.146
$ws.Range("C6").Value = 22.6
$ws.Range("C8").Value = 32

$ws.Range("C4").NumberFormat = "0.0000"
$ws.Range("C8").NumberFormat = "0.0"

$ws.Range("C4").Select()
